# Weekly CompStat (cs-en-us-pbbs) data refresh.
# Updates the report's volume/issue number, the "week covering" date range,
# and the full block of weekly/28-day/YTD/2-year crime-count + %-change
# figures for rows 14-33 (Brooklyn South precinct table), plus the
# best-fit width of column E which grew to accommodate a wider value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: issue number "18" -> "19" (last run of the rich-text A8 cell) ---
$ws.Range("A8").Value = "Volume 32   Number  19"

# --- Header: reporting week date range moves forward one week ---
$ws.Range("C9").Value = "Report Covering the Week  5/5/2025  Through  5/11/2025"

# --- Column E widened to fit a new best-fit value (closest reachable via ColumnWidth) ---
$ws.Columns.Item(5).ColumnWidth = 6.72

# --- Row 14 (Murder) ---
$ws.Range("G14").Value = 3
$ws.Range("H14").Value = -33.333333333333
$ws.Range("I14").Value = 11
$ws.Range("K14").Value = -50
$ws.Range("L14").Value = -62.068965517241
$ws.Range("M14").Value = -56
$ws.Range("N14").Value = -86.746987951807

# --- Row 15 (Rape) ---
$ws.Range("D15").Value = 4
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 14
$ws.Range("G15").Value = 15
$ws.Range("H15").Value = -6.666666666666
$ws.Range("I15").Value = 99
$ws.Range("J15").Value = 74
$ws.Range("K15").Value = 33.783783783783
$ws.Range("L15").Value = 33.783783783783
$ws.Range("M15").Value = 86.792452830188
$ws.Range("N15").Value = -47.619047619047

# --- Row 16 (Robbery) ---
$ws.Range("C16").Value = 32
$ws.Range("D16").Value = 32
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 125
$ws.Range("G16").Value = 122
$ws.Range("H16").Value = 2.459016393442
$ws.Range("I16").Value = 536
$ws.Range("J16").Value = 625
$ws.Range("K16").Value = -14.24
$ws.Range("L16").Value = -15.590551181102
$ws.Range("M16").Value = -45.361875637105
$ws.Range("N16").Value = -89.297124600639

# --- Row 17 (Fel. Assault) ---
$ws.Range("C17").Value = 84
$ws.Range("D17").Value = 68
$ws.Range("E17").Value = 23.529411764705
$ws.Range("F17").Value = 337
$ws.Range("G17").Value = 284
$ws.Range("H17").Value = 18.661971830985
$ws.Range("I17").Value = 1294
$ws.Range("J17").Value = 1203
$ws.Range("K17").Value = 7.564422277639
$ws.Range("L17").Value = 9.198312236286
$ws.Range("M17").Value = 66.323907455012
$ws.Range("N17").Value = -41.973094170403

# --- Row 18 (Burglary) ---
$ws.Range("C18").Value = 29
$ws.Range("D18").Value = 25
$ws.Range("E18").Value = 16
$ws.Range("F18").Value = 110
$ws.Range("H18").Value = 1.851851851851
$ws.Range("I18").Value = 511
$ws.Range("J18").Value = 534
$ws.Range("K18").Value = -4.307116104868
$ws.Range("L18").Value = -23.157894736842
$ws.Range("M18").Value = -57.022708158116
$ws.Range("N18").Value = -92.451994091580

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 108
$ws.Range("D19").Value = 102
$ws.Range("E19").Value = 5.882352941176
$ws.Range("F19").Value = 441
$ws.Range("G19").Value = 492
$ws.Range("H19").Value = -10.365853658536
$ws.Range("I19").Value = 1824
$ws.Range("J19").Value = 2178
$ws.Range("K19").Value = -16.253443526170
$ws.Range("L19").Value = -23.457826269408
$ws.Range("M19").Value = 1.956400223588
$ws.Range("N19").Value = -39.301164725457

# --- Row 20 (G.L.A.) ---
$ws.Range("C20").Value = 33
$ws.Range("E20").Value = -13.157894736842
$ws.Range("F20").Value = 158
$ws.Range("G20").Value = 156
$ws.Range("H20").Value = 1.282051282051
$ws.Range("I20").Value = 569
$ws.Range("J20").Value = 689
$ws.Range("K20").Value = -17.416545718432
$ws.Range("L20").Value = -6.873977086743
$ws.Range("M20").Value = -20.308123249299
$ws.Range("N20").Value = -93.470277714023

# --- Row 21 (TOTAL) ---
$ws.Range("D21").Value = 269
$ws.Range("E21").Value = 7.806691449814
$ws.Range("F21").Value = 1187
$ws.Range("G21").Value = 1180
$ws.Range("H21").Value = 0.593220338983
$ws.Range("I21").Value = 4844
$ws.Range("J21").Value = 5325
$ws.Range("K21").Value = -9.032863849765
$ws.Range("L21").Value = -13.221067717663
$ws.Range("M21").Value = -12.389220473865
$ws.Range("N21").Value = -81.368514173622

# --- Row 22 (Transit) ---
$ws.Range("C22").Value = 7
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = 40
$ws.Range("F22").Value = 18
$ws.Range("G22").Value = 16
$ws.Range("H22").Value = 12.5
$ws.Range("I22").Value = 59
$ws.Range("J22").Value = 68
$ws.Range("K22").Value = -13.235294117647
$ws.Range("L22").Value = -3.278688524590
$ws.Range("M22").Value = -37.894736842105

# --- Row 23 (Housing) ---
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = 140
$ws.Range("F23").Value = 45
$ws.Range("G23").Value = 43
$ws.Range("H23").Value = 4.651162790697
$ws.Range("I23").Value = 165
$ws.Range("J23").Value = 187
$ws.Range("K23").Value = -11.764705882352
$ws.Range("L23").Value = -4.624277456647
$ws.Range("M23").Value = 83.333333333333

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 243
$ws.Range("D24").Value = 263
$ws.Range("E24").Value = -7.604562737642
$ws.Range("F24").Value = 1020
$ws.Range("G24").Value = 1095
$ws.Range("H24").Value = -6.849315068493
$ws.Range("I24").Value = 4990
$ws.Range("J24").Value = 5244
$ws.Range("K24").Value = -4.843630816170
$ws.Range("L24").Value = -14.098812187984
$ws.Range("M24").Value = 22.214058290472

# --- Row 25 (Retail Theft) ---
$ws.Range("C25").Value = 90
$ws.Range("D25").Value = 126
$ws.Range("E25").Value = -28.571428571428
$ws.Range("F25").Value = 432
$ws.Range("G25").Value = 538
$ws.Range("H25").Value = -19.702602230483
$ws.Range("I25").Value = 2207
$ws.Range("J25").Value = 2496
$ws.Range("K25").Value = -11.578525641025
$ws.Range("L25").Value = -21.958981612447

# --- Row 26 (Misd. Assault) ---
$ws.Range("C26").Value = 94
$ws.Range("D26").Value = 140
$ws.Range("E26").Value = -32.857142857142
$ws.Range("F26").Value = 472
$ws.Range("G26").Value = 488
$ws.Range("H26").Value = -3.278688524590
$ws.Range("I26").Value = 2145
$ws.Range("J26").Value = 2137
$ws.Range("K26").Value = 0.374356574637
$ws.Range("L26").Value = 7.142857142857
$ws.Range("M26").Value = -3.854773644105

# --- Row 27 (UCR Rape*) ---
$ws.Range("C27").Value = 4
$ws.Range("D27").Value = 12
$ws.Range("E27").Value = -66.666666666666
$ws.Range("F27").Value = 17
$ws.Range("G27").Value = 26
$ws.Range("H27").Value = -34.615384615384
$ws.Range("I27").Value = 120
$ws.Range("J27").Value = 120
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 10.091743119266

# --- Row 28 (Other Sex Crimes) ---
$ws.Range("C28").Value = 14
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 53
$ws.Range("G28").Value = 59
$ws.Range("H28").Value = -10.169491525423
$ws.Range("I28").Value = 235
$ws.Range("J28").Value = 247
$ws.Range("K28").Value = -4.858299595141
$ws.Range("L28").Value = 9.302325581395

# --- Row 29 (Shooting Vic.) ---
$ws.Range("F29").Value = 13
$ws.Range("G29").Value = 7
$ws.Range("H29").Value = 85.714285714285
$ws.Range("I29").Value = 34
$ws.Range("J29").Value = 35
$ws.Range("K29").Value = -2.857142857142
$ws.Range("L29").Value = -41.379310344827
$ws.Range("M29").Value = -47.692307692307
$ws.Range("N29").Value = -87.121212121212

# --- Row 30 (Shooting Inc.) ---
$ws.Range("F30").Value = 11
$ws.Range("G30").Value = 7
$ws.Range("H30").Value = 57.142857142857
$ws.Range("I30").Value = 31
$ws.Range("J30").Value = 33
$ws.Range("K30").Value = -6.060606060606
$ws.Range("L30").Value = -36.734693877551
$ws.Range("M30").Value = -42.592592592592
$ws.Range("N30").Value = -86.808510638297

# --- Row 31 (Hate Crimes): D/E flip from "n/a" placeholders to real numbers ---
$ws.Range("D31").Value = 6
$ws.Range("E31").Value = -100
$ws.Range("F31").Value = 10
$ws.Range("G31").Value = 8
$ws.Range("H31").Value = 25
$ws.Range("I31").Value = 51
$ws.Range("J31").Value = 45
$ws.Range("K31").Value = 13.333333333333
$ws.Range("L31").Value = 96.153846153846

# --- Row 33 (Traffic Fatalities) ---
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = 3
$ws.Range("H33").Value = -66.666666666666
$ws.Range("L33").Value = -35.294117647058
